$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update the "request" progress values in column B (new user-request data) ---
# Row 3  - Celeb Profile Management
$ws.Range("B3").Value = 0.2
# Row 8  - Messaging System
$ws.Range("B8").Value = 0.3
# Row 11 - Admin Dashboard
$ws.Range("B11").Value = 0.5
# Row 12 - Content Moderation
$ws.Range("B12").Value = 0.3
# Row 16 - Forgot Password Functionality
$ws.Range("B16").Value = 1

# --- Column C cells that are manually overridden with a literal value   ---
# --- (no longer driven by the shared MROUND formula for these rows)    ---
$ws.Range("C3").Value = 0.2
$ws.Range("C11").Value = 0
$ws.Range("C12").Value = 0
$ws.Range("C16").Value = 0.8

# Row 16 is now "done" - give it the same highlighted (green) look used
# for the other completed row (row 2).
$ws.Range("B2:C2").Copy() | Out-Null
$ws.Range("B16:C16").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("B16").Value = 1
$ws.Range("C16").Value = 0.8

$excel.CutCopyMode = 0

# --- Move the active cell / selection on Sheet1 ---
$ws.Range("H19").Select() | Out-Null

$wb.Save()
